$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SeznamModulu")

# Add the two new "Kroužek 2018/2019" modules in rows 47 and 48
$ws.Range("B47").Value = "Kroužek 2018/2019 I"
$ws.Range("C47").Value2 = 43344

$ws.Range("B48").Value = "Kroužek 2018/2019 II"
$ws.Range("C48").Value2 = 43344

# Update the frozen-pane view state to match the new scroll/selection position
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$pane = $win.Panes.Item($win.Panes.Count)
$pane.ScrollRow = 28
$pane.ScrollColumn = 1
$ws.Range("C49").Select()
